$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.349435015852742
$ws.Range("D2").Value = 0.7300854328946536

$ws.Range("C3").Value = 0.3834428841291662
$ws.Range("D3").Value = 0.7050717953706505

$ws.Range("C4").Value = -1.246545215702003
$ws.Range("D4").Value = 0.2256741701118943

$ws.Range("C5").Value = -2.211383498326237
$ws.Range("D5").Value = 0.03770406171492358

$ws.Range("C6").Value = 0.06929162683616527
$ws.Range("D6").Value = 0.9453834621564245

$ws.Range("C7").Value = -1.701101503154351
$ws.Range("D7").Value = 0.1030174124732572

$ws.Range("C8").Value = -2.586232427554077
$ws.Range("D8").Value = 0.01685122040904097

$ws.Range("C9").Value = -1.738029673946147
$ws.Range("D9").Value = 0.0961849646668298

$ws.Range("C10").Value = -2.414857565502929
$ws.Range("D10").Value = 0.02449975151269879

$ws.Range("C11").Value = -1.769530255968567
$ws.Range("D11").Value = 0.09066671961494022
$ws.Range("G11").Value = "No"
